$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'69.055.50"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  +2.42%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'3.757.33"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  +1.81%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("D4").Value = "'1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'  -0.02%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'605.96"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  +1.68%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'169.59"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  +2.40%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'3.756.06"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'  +1.82%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'1.00"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'  -0.03%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'0.539"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'  +2.79%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'0.168"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'  +5.73%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'6.38"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'  +3.18%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'0.464"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  +0.55%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'38.44"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  +2.46%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'0.0000249"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  +3.93%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'4.382.99"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  +1.85%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'3.750.07"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  +1.62%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'69.077.39"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  +2.36%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'7.30"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  +2.04%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'0.114"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  +0.06%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'17.20"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  -1.97%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'10.94"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  +19.98%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'495.27"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  +0.77%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'0.730"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  +0.96%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'0.0000155"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  +12.43%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'85.58"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  +0.16%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'2.34"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'  +1.61%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = "'12.40"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  +1.96%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("D28").Value = "'10.35"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'  +3.16%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("D29").Value = "'1.00"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  -0.23%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("D30").Value = "'2.54"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  +7.79%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("D31").Value = "'3.00"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  +2.45%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("D32").Value = "'7.95"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  +4.05%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("D33").Value = "'31.97"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "'  +1.60%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("D34").Value = "'3.901.88"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "'  +1.93%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("D35").Value = "'0.109"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  +1.48%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("D36").Value = "'3.690.38"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "'  +1.69%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("D37").Value = "'1.00"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'  -0.02%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("D38").Value = "'1.02"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'  +2.29%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").Value = "'5.91"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  +2.63%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'0.133"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'  +2.17%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'0.325"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'  +0.52%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'3.05"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  +9.74%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'438.72"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  +0.61%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'48.49"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  -0.15%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'1.98"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  +2.24%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'8.48"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  +1.31%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'1.00"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  +0.01%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'40.54"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  -0.12%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = "'141.05"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  -1.30%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'0.0356"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  +2.50%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = "'2.786.91"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  +1.23%  "
$ws.Range("E51").Style = "Normal"
